$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216

$ws.Range("H32").Value = 1099.1875
$ws.Range("I32").Value = 799.5
$ws.Range("K32").Value = 799.5
$ws.Range("M32").Value = -473.5

$ws.Range("H40").Value = 5911.7646
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 6857.143
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 6857.143
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -7207.143

$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968

$ws.Range("H55").Value = 384
$ws.Range("I55").Value = 425
$ws.Range("J55").Value = 302
$ws.Range("K55").Value = 425
$ws.Range("L55").Value = 302
$ws.Range("M55").Value = -211
$ws.Range("N55").Value = -730

$ws.Range("H70").Value = 3144.5
$ws.Range("J70").Value = 3499
$ws.Range("L70").Value = 10497
$ws.Range("N70").Value = -11037

$ws.Range("H73").Value = 3144.5
$ws.Range("J73").Value = 3499
$ws.Range("L73").Value = 10497
$ws.Range("N73").Value = -12369

$ws.Range("H80").Value = 634.8946999999999
$ws.Range("J80").Value = 624.5454999999999
$ws.Range("L80").Value = 1873.6365
$ws.Range("N80").Value = -3869.6365

$ws.Range("H83").Value = 634.8946999999999
$ws.Range("J83").Value = 624.5454999999999
$ws.Range("L83").Value = 5620.9095
$ws.Range("N83").Value = -15604.9095

$ws.Range("H92").Value = 949.5
$ws.Range("I92").Value = 939.4
$ws.Range("K92").Value = 939.4
$ws.Range("M92").Value = 308.6

$ws.Range("H100").Value = 2197.9
$ws.Range("I100").Value = 2829.6667
$ws.Range("J100").Value = 1250.25
$ws.Range("K100").Value = 2829.6667
$ws.Range("L100").Value = 1250.25
$ws.Range("M100").Value = -2288.6667
$ws.Range("N100").Value = -2332.25

$ws.Range("H107").Value = 738
$ws.Range("I107").Value = 721.1429000000001
$ws.Range("J107").Value = 752.75
$ws.Range("K107").Value = 721.1429000000001
$ws.Range("L107").Value = 752.75
$ws.Range("M107").Value = 1198.8571
$ws.Range("N107").Value = -4592.75

$ws.Range("H115").Value = 117.5
$ws.Range("I115").Value = 117.5
$ws.Range("K115").Value = 352.5
$ws.Range("M115").Value = 1214.5

$ws.Range("H137").Value = 1507.6666
$ws.Range("I137").Value = 1507.6666
$ws.Range("K137").Value = 4522.9998
$ws.Range("M137").Value = -1972.9998

$ws.Range("H138").Value = 1971.619
$ws.Range("J138").Value = 2913.6365
$ws.Range("L138").Value = 8740.9095
$ws.Range("N138").Value = -19020.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 763.125
$ws.Range("I97").Value = 733.25
$ws.Range("J97").Value = 852.75
$ws.Range("K97").Value = 733.25
$ws.Range("L97").Value = 852.75
$ws.Range("M97").Value = -237.25
$ws.Range("N97").Value = -1844.75

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H86").Value = 4353.857
$ws.Range("I86").Value = 4619.25
$ws.Range("K86").Value = 4619.25
$ws.Range("M86").Value = -3496.25

$ws.Range("H89").Value = 4353.857
$ws.Range("I89").Value = 4619.25
$ws.Range("K89").Value = 23096.25
$ws.Range("M89").Value = -17480.25

$ws.Range("H94").Value = 1046.7273
$ws.Range("I94").Value = 991.4
$ws.Range("J94").Value = 1600
$ws.Range("K94").Value = 991.4
$ws.Range("L94").Value = 1600
$ws.Range("M94").Value = -540.4
$ws.Range("N94").Value = -2502

$ws.Range("H99").Value = 4199.8887
$ws.Range("I99").Value = 4237.5
$ws.Range("J99").Value = 3899
$ws.Range("K99").Value = 4237.5
$ws.Range("L99").Value = 3899
$ws.Range("M99").Value = -2739.5
$ws.Range("N99").Value = -6895

$ws.Range("H107").Value = 1152.6923
$ws.Range("I107").Value = 1031.7778
$ws.Range("J107").Value = 1424.75
$ws.Range("K107").Value = 1031.7778
$ws.Range("L107").Value = 1424.75
$ws.Range("M107").Value = 888.2221999999999
$ws.Range("N107").Value = -5264.75

$ws.Range("H134").Value = 8564.223
$ws.Range("I134").Value = 8564.223
$ws.Range("K134").Value = 25692.669
$ws.Range("M134").Value = -23157.669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 72.5
$ws.Range("K7").Value = 72.5
$ws.Range("M7").Value = 40.5

$ws.Range("H22").Value = 753.25
$ws.Range("I22").Value = 381.5
$ws.Range("K22").Value = 381.5
$ws.Range("M22").Value = -31.5

$ws.Range("H107").Value = 678.5
$ws.Range("I107").Value = 657.8333
$ws.Range("J107").Value = 719.8333
$ws.Range("K107").Value = 657.8333
$ws.Range("L107").Value = 719.8333
$ws.Range("M107").Value = 1262.1667
$ws.Range("N107").Value = -4559.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2230
$ws.Range("I4").Value = 2475
$ws.Range("K4").Value = 7425
$ws.Range("M4").Value = -7313

$ws.Range("H11").Value = 500075
$ws.Range("I11").Value = 666700
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 2000100
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -1999960
$ws.Range("N11").Value = -880

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H47").Value = 852.6
$ws.Range("J47").Value = 916.6667
$ws.Range("L47").Value = 2750.0001
$ws.Range("N47").Value = -3612.0001

$ws.Range("H75").Value = 4800
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 4800
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H119").Value = 200
$ws.Range("I119").Value = 200
$ws.Range("K119").Value = 600
$ws.Range("M119").Value = 4238

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 39000
$ws.Range("J18").Value = 39000
$ws.Range("L18").Value = 39000
$ws.Range("N18").Value = -39586

$ws.Range("H97").Value = 561.5
$ws.Range("I97").Value = 561.5
$ws.Range("K97").Value = 561.5
$ws.Range("M97").Value = -65.5

$ws.Range("H107").Value = 186.4
$ws.Range("I107").Value = 220.14285
$ws.Range("J107").Value = 107.666664
$ws.Range("K107").Value = 220.14285
$ws.Range("L107").Value = 107.666664
$ws.Range("M107").Value = 1699.85715
$ws.Range("N107").Value = -3947.666664

$ws.Range("H126").Value = 14822.2
$ws.Range("I126").Value = 11027.75
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 33083.25
$ws.Range("L126").Value = 90000
$ws.Range("M126").Value = -30613.25
$ws.Range("N126").Value = -94940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2424.75
$ws.Range("I16").Value = 2349.5
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2349.5
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -2179.5
$ws.Range("N16").Value = -2840

$ws.Range("H20").Value = 999
$ws.Range("I20").Value = 999
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 999
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -773
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 1832.8334
$ws.Range("I22").Value = 1666
$ws.Range("K22").Value = 1666
$ws.Range("M22").Value = -1371

$ws.Range("H27").Value = 1832.8334
$ws.Range("I27").Value = 1666
$ws.Range("K27").Value = 1666
$ws.Range("M27").Value = -1559

$ws.Range("H39").Value = 21721.25
$ws.Range("I39").Value = 17554.5
$ws.Range("K39").Value = 17554.5
$ws.Range("M39").Value = -17094.5

$ws.Range("H40").Value = 7749.75
$ws.Range("I40").Value = 6999.6665
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 6999.6665
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -6863.6665
$ws.Range("N40").Value = -10272

$ws.Range("H43").Value = 10500
$ws.Range("J43").Value = 10500
$ws.Range("L43").Value = 10500
$ws.Range("N43").Value = -10886

$ws.Range("H46").Value = 3899.9
$ws.Range("I46").Value = 2777.3333
$ws.Range("K46").Value = 2777.3333
$ws.Range("M46").Value = -2589.3333

$ws.Range("H55").Value = 211
$ws.Range("I55").Value = 182.16667
$ws.Range("J55").Value = 239.83333
$ws.Range("K55").Value = 182.16667
$ws.Range("L55").Value = 239.83333
$ws.Range("M55").Value = -9.166670000000011
$ws.Range("N55").Value = -585.8333299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H111").Value = 21000
$ws.Range("J111").Value = 21000
$ws.Range("L111").Value = 21000
$ws.Range("N111").Value = -29180
